# testdata.xlsx update - 15 march 2020
# 1. Insert a new "test_suite" sheet as the first sheet, driving which
#    test-case sheets should run.
# 2. Fix the "lastanme" header typo on AddCustomerTest (-> "lastname"),
#    add a "runmode" column, and append three more test rows.
# 3. Update OpenAccountTest's customer name to "Raman Arora".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "test_suite" sheet and move it to the front
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "test_suite"
$newSheet.Move($wb.Worksheets.Item(1))

$suite = $wb.Worksheets.Item("test_suite")

$suite.Cells.Item(1,1).Value = "TCID"
$suite.Cells.Item(1,2).Value = "Runmode"
$suite.Cells.Item(2,1).Value = "BankManagerLoginTest"
$suite.Cells.Item(2,2).Value = "Y"
$suite.Cells.Item(3,1).Value = "AddCustomerTest"
$suite.Cells.Item(3,2).Value = "Y"
$suite.Cells.Item(4,1).Value = "OpenAccountTest"
$suite.Cells.Item(4,2).Value = "Y"

$suite.Columns.Item(1).AutoFit() | Out-Null
$suite.Range("C7").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) AddCustomerTest: fix header, add runmode column + new rows
# ---------------------------------------------------------------------
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")

$addCustomer.Cells.Item(1,2).Value = "lastname"
$addCustomer.Cells.Item(1,5).Value = "runmode"

$addCustomer.Cells.Item(2,5).Value = "Y"

$addCustomer.Cells.Item(3,1).Value = "Vishal"
$addCustomer.Cells.Item(3,2).Value = "Arora"
$addCustomer.Cells.Item(3,3).Value = 35435
$addCustomer.Cells.Item(3,4).Value = "Customer added successfully"
$addCustomer.Cells.Item(3,5).Value = "N"

$addCustomer.Cells.Item(4,1).Value = "Soniya"
$addCustomer.Cells.Item(4,2).Value = "Arora"
$addCustomer.Cells.Item(4,3).Value = 35435
$addCustomer.Cells.Item(4,4).Value = "Customer added successfully"
$addCustomer.Cells.Item(4,5).Value = "Y"

$addCustomer.Cells.Item(5,1).Value = "Rohit"
$addCustomer.Cells.Item(5,2).Value = "Sehgal"
$addCustomer.Cells.Item(5,3).Value = 35435
$addCustomer.Cells.Item(5,4).Value = "Customer added successfully"
$addCustomer.Cells.Item(5,5).Value = "N"

$addCustomer.Activate()
$addCustomer.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) OpenAccountTest: customer becomes "Raman Arora"
# ---------------------------------------------------------------------
$openAccount = $wb.Worksheets.Item("OpenAccountTest")
$openAccount.Cells.Item(2,1).Value = "Raman Arora"
$openAccount.Columns.Item(1).AutoFit() | Out-Null
$openAccount.Range("A2").Select() | Out-Null
